# Remove the "biosat" and "O2_Ar_ratio" attribute rows (rows 7 and 8)
# from the ColumnHeaders sheet, per commit message:
# "removed O2Ar & biosat from ncp output, wrote csvs for qc"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 7 (biosat) and 8 (O2_Ar_ratio); rows below shift up,
# so former row 9 (ncp) becomes row 7 and former row 10 (k) becomes row 8.
$ws.Rows("7:8").Delete()

# Match the selection left behind by the edit in Excel.
$ws.Range("A7:XFD8").Select()
